# edit.ps1 — applies the cryptos.xlsx price/volume refresh described in the commit
# "Updated cryptos list on Wed Oct 25 05:27:29 UTC 2023 with GitHub Actions".
#
# Columns: A=rank(unchanged) B=Coin C=Link D=Price(text) E=Volume(1h)(text, "  +x.xx%  ")
# D/E are stored as plain text in the source feed, so any D value that LOOKS like a
# normal decimal number (single "." , e.g. "227.66") would otherwise get silently
# auto-coerced to a numeric cell by Excel on assignment. Forcing NumberFormat "@"
# (Text) on just those target cells before writing keeps them as text, matching the
# original inlineStr typing. Multi-dot values (e.g. "34.058.58") and the padded
# "  +n.nn%  " strings are never mistaken for numbers, so no format coercion is needed there.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "34.058.58"
$ws.Cells.Item(2, 5).Value = "  -1.19%  "

$ws.Cells.Item(3, 4).Value = "1.788.89"
$ws.Cells.Item(3, 5).Value = "  -2.27%  "

$ws.Cells.Item(4, 5).Value = "  +0.20%  "

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "227.66"
$ws.Cells.Item(5, 5).Value = "  -1.56%  "

$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "0.553"
$ws.Cells.Item(6, 5).Value = "  +1.45%  "

$ws.Cells.Item(7, 5).Value = "  +0.18%  "

$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "31.51"
$ws.Cells.Item(8, 5).Value = "  -0.39%  "

$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "46.46"
$ws.Cells.Item(9, 5).Value = "  +1.52%  "

$ws.Cells.Item(10, 5).Value = "  -0.63%  "

$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "0.0660"
$ws.Cells.Item(11, 5).Value = "  -2.73%  "

$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "0.0928"
$ws.Cells.Item(12, 5).Value = "  -0.49%  "

$ws.Cells.Item(13, 4).Value = "2.049.03"
$ws.Cells.Item(13, 5).Value = "  -1.91%  "

$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "11.53"
$ws.Cells.Item(14, 5).Value = "  +11.57%  "

$ws.Cells.Item(15, 4).Value = "1.790.84"
$ws.Cells.Item(15, 5).Value = "  -1.94%  "

$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "0.636"
$ws.Cells.Item(16, 5).Value = "  -2.03%  "

$ws.Cells.Item(17, 4).Value = "34.075.88"
$ws.Cells.Item(17, 5).Value = "  -1.09%  "

$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "4.23"
$ws.Cells.Item(18, 5).Value = "  -3.25%  "

$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "69.48"
$ws.Cells.Item(19, 5).Value = "  -1.07%  "

$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "253.17"
$ws.Cells.Item(20, 5).Value = "  -2.94%  "

$ws.Cells.Item(21, 4).Value = "0.0₃0742"
$ws.Cells.Item(21, 5).Value = "  -1.43%  "

$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "0.999"
$ws.Cells.Item(22, 5).Value = "  +0.00%  "

$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "10.49"
$ws.Cells.Item(23, 5).Value = "  -0.61%  "

$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "4.28"
$ws.Cells.Item(24, 5).Value = "  -2.77%  "

$ws.Cells.Item(25, 5).Value = "  -2.58%  "

$ws.Cells.Item(26, 5).Value = "  -2.14%  "

$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "16.59"
$ws.Cells.Item(27, 5).Value = "  -1.55%  "

$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "7.04"
$ws.Cells.Item(28, 5).Value = "  -2.42%  "

$ws.Cells.Item(29, 5).Value = "  -1.95%  "

$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "1.00"
$ws.Cells.Item(30, 5).Value = "  +0.24%  "

$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "3.82"
$ws.Cells.Item(31, 5).Value = "  -1.05%  "

$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "0.0516"
$ws.Cells.Item(32, 5).Value = "  -0.36%  "

$ws.Cells.Item(33, 5).Value = "  -1.24%  "

$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "3.61"
$ws.Cells.Item(34, 5).Value = "  +0.70%  "

$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "1.85"
$ws.Cells.Item(35, 5).Value = "  -0.21%  "

$ws.Cells.Item(36, 4).Value = "1.454.80"
$ws.Cells.Item(36, 5).Value = "  -8.38%  "

$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "1.07"
$ws.Cells.Item(37, 5).Value = "  -0.09%  "

$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "0.630"
$ws.Cells.Item(38, 5).Value = "  -0.69%  "

$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "0.0187"
$ws.Cells.Item(39, 5).Value = "  -1.77%  "

$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "83.45"
$ws.Cells.Item(40, 5).Value = "  -2.47%  "

$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "2.83"
$ws.Cells.Item(41, 5).Value = "  -1.64%  "

$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "2.35"
$ws.Cells.Item(42, 5).Value = "  -0.46%  "

$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "0.904"
$ws.Cells.Item(43, 5).Value = "  -1.91%  "

$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "2.08"
$ws.Cells.Item(44, 5).Value = "  -3.59%  "

$ws.Cells.Item(45, 5).Value = "  -1.84%  "

$ws.Cells.Item(46, 5).Value = "  +0.60%  "

$ws.Cells.Item(47, 4).Value = "1.946.13"
$ws.Cells.Item(47, 5).Value = "  -1.68%  "

$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "5.78"
$ws.Cells.Item(48, 5).Value = "  +0.32%  "

$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "51.30"
$ws.Cells.Item(51, 5).Value = "  -3.78%  "

# Rows 49/50: the feed re-sorted two coins with near-identical rank; InjectiveProtocol
# and PaxDollar swap rows, each carrying its own refreshed Price/Volume values.
$ws.Cells.Item(49, 2).Value = "PaxDollar"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "1.00"
$ws.Cells.Item(49, 5).Value = "  +0.16%  "

$ws.Cells.Item(50, 2).Value = "InjectiveProtocol"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "11.98"
$ws.Cells.Item(50, 5).Value = "  +5.79%  "

